$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-04 Wednesday", "2024-12-05 Thursday"),
    @("871×8=6968", "718×5=3590"),
    @("576×6=3456", "519×7=3633"),
    @("296×9=2664", "618×5=3090"),
    @("925×6=5550", "137×5=685"),
    @("863×5=4315", "307×7=2149"),
    @("348×2=696", "676×5=3380"),
    @("442×6=2652", "689×3=2067"),
    @("242×5=1210", "897×8=7176"),
    @("592×7=4144", "157×5=785"),
    @("990×6=5940", "901×8=7208"),
    @("277×2=554", "792×9=7128"),
    @("438×2=876", "809×7=5663"),
    @("781×7=5467", "286×7=2002"),
    @("838×5=4190", "495×6=2970"),
    @("370×8=2960", "605×6=3630"),
    @("405×9=3645", "858×2=1716"),
    @("454×6=2724", "232×7=1624"),
    @("417×9=3753", "455×8=3640"),
    @("847×3=2541", "946×8=7568"),
    @("633×4=2532", "342×2=684"),
    @("300×5=1500", "251×6=1506"),
    @("779×5=3895", "679×2=1358"),
    @("874×8=6992", "506×8=4048"),
    @("152×5=760", "684×3=2052"),
    @("849×6=5094", "218×3=654")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
